$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 16.71895933333333
$ws.Cells.Item(2, 8).Value = 50.156878
$ws.Cells.Item(2, 9).Value = 0.02912144738161902
$ws.Cells.Item(2, 10).Value = 0.03059269312988411
$ws.Cells.Item(2, 13).Value = 12.431794
$ws.Cells.Item(2, 14).Value = 37.295382
$ws.Cells.Item(2, 15).Value = 0.6267040910788743
$ws.Cells.Item(2, 16).Value = 0.7075740515758999
$ws.Cells.Item(2, 17).Value = 207.8466583263774
$ws.Cells.Item(2, 18).Value = 1870.619924937396
$ws.Cells.Item(2, 19).Value = 0.01825053021219881
$ws.Cells.Item(2, 20).Value = 0.0216465958265303

$ws.Cells.Item(3, 7).Value = 16.71895933333333
$ws.Cells.Item(3, 8).Value = 50.156878
$ws.Cells.Item(3, 9).Value = 0.02912144738161902
$ws.Cells.Item(3, 10).Value = 0.03059269312988411
$ws.Cells.Item(3, 15).Value = 0.0264162940991436
$ws.Cells.Item(3, 16).Value = 0.0298250554119953
$ws.Cells.Item(3, 17).Value = 8.760974329083776
$ws.Cells.Item(3, 18).Value = 78.84876896175399
$ws.Cells.Item(3, 19).Value = 0.0007692807186255832
$ws.Cells.Item(3, 20).Value = 0.0009124287678009617

$ws.Cells.Item(4, 7).Value = 16.71895933333333
$ws.Cells.Item(4, 8).Value = 50.156878
$ws.Cells.Item(4, 9).Value = 0.02912144738161902
$ws.Cells.Item(4, 10).Value = 0.03059269312988411
$ws.Cells.Item(4, 13).Value = 0.03915333333333333
$ws.Cells.Item(4, 14).Value = 0.11746
$ws.Cells.Item(4, 15).Value = 0.001973774193762771
$ws.Cells.Item(4, 16).Value = 0.002228470219130754
$ws.Cells.Item(4, 17).Value = 0.6546029877644445
$ws.Cells.Item(4, 18).Value = 5.89142688988
$ws.Cells.Item(4, 19).Value = 0.00005747916132686002
$ws.Cells.Item(4, 20).Value = 0.00006817490556295276

$ws.Cells.Item(5, 7).Value = 16.71895933333333
$ws.Cells.Item(5, 8).Value = 50.156878
$ws.Cells.Item(5, 9).Value = 0.02912144738161902
$ws.Cells.Item(5, 10).Value = 0.03059269312988411
$ws.Cells.Item(5, 13).Value = 6.8015495
$ws.Cells.Item(5, 14).Value = 13.603099
$ws.Cells.Item(5, 15).Value = 0.3428756056708687
$ws.Cells.Item(5, 16).Value = 0.2580802061075034
$ws.Cells.Item(5, 17).Value = 113.7148294941537
$ws.Cells.Item(5, 18).Value = 682.288976964922
$ws.Cells.Item(5, 19).Value = 0.009985033908984952
$ws.Cells.Item(5, 20).Value = 0.007895368548344094

$ws.Cells.Item(6, 7).Value = 16.71895933333333
$ws.Cells.Item(6, 8).Value = 50.156878
$ws.Cells.Item(6, 9).Value = 0.02912144738161902
$ws.Cells.Item(6, 10).Value = 0.03059269312988411
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.04027333333333333
$ws.Cells.Item(6, 14).Value = 0.12082
$ws.Cells.Item(6, 15).Value = 0.002030234957350741
$ws.Cells.Item(6, 16).Value = 0.002292216685470609
$ws.Cells.Item(6, 17).Value = 0.6733282222177779
$ws.Cells.Item(6, 18).Value = 6.059953999959999
$ws.Cells.Item(6, 19).Value = 0.00005912338048281313
$ws.Cells.Item(6, 20).Value = 0.00007012508164580243

$ws.Cells.Item(7, 9).Value = 0.2708539632042961
$ws.Cells.Item(7, 10).Value = 0.2845377865576845
$ws.Cells.Item(7, 13).Value = 12.431794
$ws.Cells.Item(7, 14).Value = 37.295382
$ws.Cells.Item(7, 15).Value = 0.6267040910788743
$ws.Cells.Item(7, 16).Value = 0.7075740515758999
$ws.Cells.Item(7, 17).Value = 1933.148802967866
$ws.Cells.Item(7, 18).Value = 17398.3392267108
$ws.Cells.Item(7, 19).Value = 0.1697452868250592
$ws.Cells.Item(7, 20).Value = 0.2013315544610595

$ws.Cells.Item(8, 9).Value = 0.2708539632042961
$ws.Cells.Item(8, 10).Value = 0.2845377865576845
$ws.Cells.Item(8, 15).Value = 0.0264162940991436
$ws.Cells.Item(8, 16).Value = 0.0298250554119953
$ws.Cells.Item(8, 19).Value = 0.007154957949923303
$ws.Cells.Item(8, 20).Value = 0.008486355250889433

$ws.Cells.Item(9, 9).Value = 0.2708539632042961
$ws.Cells.Item(9, 10).Value = 0.2845377865576845
$ws.Cells.Item(9, 13).Value = 0.03915333333333333
$ws.Cells.Item(9, 14).Value = 0.11746
$ws.Cells.Item(9, 15).Value = 0.001973774193762771
$ws.Cells.Item(9, 16).Value = 0.002228470219130754
$ws.Cells.Item(9, 17).Value = 6.08835856398
$ws.Cells.Item(9, 18).Value = 54.79522707582
$ws.Cells.Item(9, 19).Value = 0.0005346045628510106
$ws.Cells.Item(9, 20).Value = 0.0006340839835611829

$ws.Cells.Item(10, 9).Value = 0.2708539632042961
$ws.Cells.Item(10, 10).Value = 0.2845377865576845
$ws.Cells.Item(10, 13).Value = 6.8015495
$ws.Cells.Item(10, 14).Value = 13.603099
$ws.Cells.Item(10, 15).Value = 0.3428756056708687
$ws.Cells.Item(10, 16).Value = 0.2580802061075034
$ws.Cells.Item(10, 17).Value = 1057.643593052756
$ws.Cells.Item(10, 18).Value = 6345.861558316533
$ws.Cells.Item(10, 19).Value = 0.09286921668202819
$ws.Cells.Item(10, 20).Value = 0.07343357060018002

$ws.Cells.Item(11, 9).Value = 0.2708539632042961
$ws.Cells.Item(11, 10).Value = 0.2845377865576845
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.04027333333333333
$ws.Cells.Item(11, 14).Value = 0.12082
$ws.Cells.Item(11, 15).Value = 0.002030234957350741
$ws.Cells.Item(11, 16).Value = 0.002292216685470609
$ws.Cells.Item(11, 17).Value = 6.262518999660001
$ws.Cells.Item(11, 18).Value = 56.36267099694
$ws.Cells.Item(11, 19).Value = 0.0005498971844343531
$ws.Cells.Item(11, 20).Value = 0.0006522222619943993

$ws.Cells.Item(12, 7).Value = 194.8548433333333
$ws.Cells.Item(12, 8).Value = 584.56453
$ws.Cells.Item(12, 9).Value = 0.3394024086099587
$ws.Cells.Item(12, 10).Value = 0.3565493705749576
$ws.Cells.Item(12, 13).Value = 12.431794
$ws.Cells.Item(12, 14).Value = 37.295382
$ws.Cells.Item(12, 15).Value = 0.6267040910788743
$ws.Cells.Item(12, 16).Value = 0.7075740515758999
$ws.Cells.Item(12, 17).Value = 2422.395272222273
$ws.Cells.Item(12, 18).Value = 21801.55745000046
$ws.Cells.Item(12, 19).Value = 0.2127048779978849
$ws.Cells.Item(12, 20).Value = 0.2522850827245597

$ws.Cells.Item(13, 7).Value = 194.8548433333333
$ws.Cells.Item(13, 8).Value = 584.56453
$ws.Cells.Item(13, 9).Value = 0.3394024086099587
$ws.Cells.Item(13, 10).Value = 0.3565493705749576
$ws.Cells.Item(13, 15).Value = 0.0264162940991436
$ws.Cells.Item(13, 16).Value = 0.0298250554119953
$ws.Cells.Item(13, 17).Value = 102.1067308260878
$ws.Cells.Item(13, 18).Value = 918.9605774347899
$ws.Cells.Item(13, 19).Value = 0.008965753843798378
$ws.Cells.Item(13, 20).Value = 0.01063410473451016

$ws.Cells.Item(14, 7).Value = 194.8548433333333
$ws.Cells.Item(14, 8).Value = 584.56453
$ws.Cells.Item(14, 9).Value = 0.3394024086099587
$ws.Cells.Item(14, 10).Value = 0.3565493705749576
$ws.Cells.Item(14, 13).Value = 0.03915333333333333
$ws.Cells.Item(14, 14).Value = 0.11746
$ws.Cells.Item(14, 15).Value = 0.001973774193762771
$ws.Cells.Item(14, 16).Value = 0.002228470219130754
$ws.Cells.Item(14, 17).Value = 7.629216632644444
$ws.Cells.Item(14, 18).Value = 68.6629496938
$ws.Cells.Item(14, 19).Value = 0.0006699037154152638
$ws.Cells.Item(14, 20).Value = 0.0007945596539761081

$ws.Cells.Item(15, 7).Value = 194.8548433333333
$ws.Cells.Item(15, 8).Value = 584.56453
$ws.Cells.Item(15, 9).Value = 0.3394024086099587
$ws.Cells.Item(15, 10).Value = 0.3565493705749576
$ws.Cells.Item(15, 13).Value = 6.8015495
$ws.Cells.Item(15, 14).Value = 13.603099
$ws.Cells.Item(15, 15).Value = 0.3428756056708687
$ws.Cells.Item(15, 16).Value = 0.2580802061075034
$ws.Cells.Item(15, 17).Value = 1325.314862246412
$ws.Cells.Item(15, 18).Value = 7951.88917347847
$ws.Cells.Item(15, 19).Value = 0.1163728064182913
$ws.Cells.Item(15, 20).Value = 0.09201833504548565

$ws.Cells.Item(16, 7).Value = 194.8548433333333
$ws.Cells.Item(16, 8).Value = 584.56453
$ws.Cells.Item(16, 9).Value = 0.3394024086099587
$ws.Cells.Item(16, 10).Value = 0.3565493705749576
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.04027333333333333
$ws.Cells.Item(16, 14).Value = 0.12082
$ws.Cells.Item(16, 15).Value = 0.002030234957350741
$ws.Cells.Item(16, 16).Value = 0.002292216685470609
$ws.Cells.Item(16, 17).Value = 7.847454057177777
$ws.Cells.Item(16, 18).Value = 70.62708651459999
$ws.Cells.Item(16, 19).Value = 0.0006890666345689783
$ws.Cells.Item(16, 20).Value = 0.0008172884164259612

$ws.Cells.Item(17, 7).Value = 82.82950199999999
$ws.Cells.Item(17, 8).Value = 165.659004
$ws.Cells.Item(17, 9).Value = 0.1442742299952585
$ws.Cells.Item(17, 10).Value = 0.1010420758958371
$ws.Cells.Item(17, 13).Value = 12.431794
$ws.Cells.Item(17, 14).Value = 37.295382
$ws.Cells.Item(17, 15).Value = 0.6267040910788743
$ws.Cells.Item(17, 16).Value = 0.7075740515758999
$ws.Cells.Item(17, 17).Value = 1029.719305986588
$ws.Cells.Item(17, 18).Value = 6178.315835919528
$ws.Cells.Item(17, 19).Value = 0.09041725017528293
$ws.Cells.Item(17, 20).Value = 0.07149475102125706

$ws.Cells.Item(18, 7).Value = 82.82950199999999
$ws.Cells.Item(18, 8).Value = 165.659004
$ws.Cells.Item(18, 9).Value = 0.1442742299952585
$ws.Cells.Item(18, 10).Value = 0.1010420758958371
$ws.Cells.Item(18, 15).Value = 0.0264162940991436
$ws.Cells.Item(18, 16).Value = 0.0298250554119953
$ws.Cells.Item(18, 17).Value = 43.40384627086199
$ws.Cells.Item(18, 18).Value = 260.423077625172
$ws.Cells.Item(18, 19).Value = 0.003811190490482233
$ws.Cells.Item(18, 20).Value = 0.003013585512536378

$ws.Cells.Item(19, 7).Value = 82.82950199999999
$ws.Cells.Item(19, 8).Value = 165.659004
$ws.Cells.Item(19, 9).Value = 0.1442742299952585
$ws.Cells.Item(19, 10).Value = 0.1010420758958371
$ws.Cells.Item(19, 13).Value = 0.03915333333333333
$ws.Cells.Item(19, 14).Value = 0.11746
$ws.Cells.Item(19, 15).Value = 0.001973774193762771
$ws.Cells.Item(19, 16).Value = 0.002228470219130754
$ws.Cells.Item(19, 17).Value = 3.243051101639999
$ws.Cells.Item(19, 18).Value = 19.45830660984
$ws.Cells.Item(19, 19).Value = 0.0002847647519896359
$ws.Cells.Item(19, 20).Value = 0.0002251692570130225

$ws.Cells.Item(20, 7).Value = 82.82950199999999
$ws.Cells.Item(20, 8).Value = 165.659004
$ws.Cells.Item(20, 9).Value = 0.1442742299952585
$ws.Cells.Item(20, 10).Value = 0.1010420758958371
$ws.Cells.Item(20, 13).Value = 6.8015495
$ws.Cells.Item(20, 14).Value = 13.603099
$ws.Cells.Item(20, 15).Value = 0.3428756056708687
$ws.Cells.Item(20, 16).Value = 0.2580802061075034
$ws.Cells.Item(20, 17).Value = 563.3689579133489
$ws.Cells.Item(20, 18).Value = 2253.475831653396
$ws.Cells.Item(20, 19).Value = 0.04946811399232246
$ws.Cells.Item(20, 20).Value = 0.02607695977272765

$ws.Cells.Item(21, 7).Value = 82.82950199999999
$ws.Cells.Item(21, 8).Value = 165.659004
$ws.Cells.Item(21, 9).Value = 0.1442742299952585
$ws.Cells.Item(21, 10).Value = 0.1010420758958371
$ws.Cells.Item(21, 11).Value = 1
$ws.Cells.Item(21, 12).Value = 0.3333333333333333
$ws.Cells.Item(21, 13).Value = 0.04027333333333333
$ws.Cells.Item(21, 14).Value = 0.12082
$ws.Cells.Item(21, 15).Value = 0.002030234957350741
$ws.Cells.Item(21, 16).Value = 0.002292216685470609
$ws.Cells.Item(21, 17).Value = 3.33582014388
$ws.Cells.Item(21, 18).Value = 20.01492086328
$ws.Cells.Item(21, 19).Value = 0.0002929105851812346
$ws.Cells.Item(21, 20).Value = 0.0002316103323030256

$ws.Cells.Item(22, 7).Value = 124.2078576666667
$ws.Cells.Item(22, 8).Value = 372.623573
$ws.Cells.Item(22, 9).Value = 0.2163479508088675
$ws.Cells.Item(22, 10).Value = 0.2272780738416368
$ws.Cells.Item(22, 13).Value = 12.431794
$ws.Cells.Item(22, 14).Value = 37.295382
$ws.Cells.Item(22, 15).Value = 0.6267040910788743
$ws.Cells.Item(22, 16).Value = 0.7075740515758999
$ws.Cells.Item(22, 17).Value = 1544.126499693321
$ws.Cells.Item(22, 18).Value = 13897.13849723989
$ws.Cells.Item(22, 19).Value = 0.1355861458684483
$ws.Cells.Item(22, 20).Value = 0.1608160675424935

$ws.Cells.Item(23, 7).Value = 124.2078576666667
$ws.Cells.Item(23, 8).Value = 372.623573
$ws.Cells.Item(23, 9).Value = 0.2163479508088675
$ws.Cells.Item(23, 10).Value = 0.2272780738416368
$ws.Cells.Item(23, 15).Value = 0.0264162940991436
$ws.Cells.Item(23, 16).Value = 0.0298250554119953
$ws.Cells.Item(23, 17).Value = 65.08669772995987
$ws.Cells.Item(23, 18).Value = 585.7802795696389
$ws.Cells.Item(23, 19).Value = 0.005715111096314097
$ws.Cells.Item(23, 20).Value = 0.006778581146258378

$ws.Cells.Item(24, 7).Value = 124.2078576666667
$ws.Cells.Item(24, 8).Value = 372.623573
$ws.Cells.Item(24, 9).Value = 0.2163479508088675
$ws.Cells.Item(24, 10).Value = 0.2272780738416368
$ws.Cells.Item(24, 13).Value = 0.03915333333333333
$ws.Cells.Item(24, 14).Value = 0.11746
$ws.Cells.Item(24, 15).Value = 0.001973774193762771
$ws.Cells.Item(24, 16).Value = 0.002228470219130754
$ws.Cells.Item(24, 17).Value = 4.863151653842221
$ws.Cells.Item(24, 18).Value = 43.76836488458
$ws.Cells.Item(24, 19).Value = 0.0004270220021800001
$ws.Cells.Item(24, 20).Value = 0.0005064824190174882

$ws.Cells.Item(25, 7).Value = 124.2078576666667
$ws.Cells.Item(25, 8).Value = 372.623573
$ws.Cells.Item(25, 9).Value = 0.2163479508088675
$ws.Cells.Item(25, 10).Value = 0.2272780738416368
$ws.Cells.Item(25, 13).Value = 6.8015495
$ws.Cells.Item(25, 14).Value = 13.603099
$ws.Cells.Item(25, 15).Value = 0.3428756056708687
$ws.Cells.Item(25, 16).Value = 0.2580802061075034
$ws.Cells.Item(25, 17).Value = 844.8058922087878
$ws.Cells.Item(25, 18).Value = 5068.835353252726
$ws.Cells.Item(25, 19).Value = 0.07418043466924176
$ws.Cells.Item(25, 20).Value = 0.058655972140766

$ws.Cells.Item(26, 7).Value = 124.2078576666667
$ws.Cells.Item(26, 8).Value = 372.623573
$ws.Cells.Item(26, 9).Value = 0.2163479508088675
$ws.Cells.Item(26, 10).Value = 0.2272780738416368
$ws.Cells.Item(26, 11).Value = 1
$ws.Cells.Item(26, 12).Value = 0.3333333333333333
$ws.Cells.Item(26, 13).Value = 0.04027333333333333
$ws.Cells.Item(26, 14).Value = 0.12082
$ws.Cells.Item(26, 15).Value = 0.002030234957350741
$ws.Cells.Item(26, 16).Value = 0.002292216685470609
$ws.Cells.Item(26, 17).Value = 5.0022644544288895
$ws.Cells.Item(26, 18).Value = 45.02038008986
$ws.Cells.Item(26, 19).Value = 0.0004392371726833613
$ws.Cells.Item(26, 20).Value = 0.000520970593101421
